$d = $word.ActiveDocument

$replacements = @(
    @("2024-02-08 Thursday", "2024-02-09 Friday"),
    @("327×8=", "874×3="),
    @("618×3=", "722×9="),
    @("334×9=", "256×4="),
    @("118×7=", "590×5="),
    @("401×8=", "587×8="),
    @("121×3=", "245×7="),
    @("988×9=", "542×4="),
    @("492×9=", "778×8="),
    @("115×4=", "435×7="),
    @("765×8=", "878×7="),
    @("370×6=", "625×3="),
    @("126×2=", "581×4="),
    @("449×9=", "452×7="),
    @("613×4=", "552×5="),
    @("673×8=", "247×7="),
    @("443×5=", "851×4="),
    @("273×2=", "154×8="),
    @("468×7=", "159×5="),
    @("195×2=", "187×6="),
    @("223×7=", "421×3="),
    @("810×9=", "252×4="),
    @("834×8=", "124×2="),
    @("609×6=", "955×4="),
    @("103×9=", "213×3="),
    @("450×2=", "277×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
